# Auto-generated COM-interop script applying the Lamia_Profits.xlsx commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1502.4667
$ws.Range("I33").Value = 278.7
$ws.Range("J33").Value = 3950
$ws.Range("K33").Value = 278.7
$ws.Range("L33").Value = 3950
$ws.Range("M33").Value = -49.69999999999999
$ws.Range("N33").Value = -4408

# Row 50
$ws.Range("H50").Value = 372.75
$ws.Range("J50").Value = 372.75
$ws.Range("L50").Value = 1118.25
$ws.Range("N50").Value = -2068.25

# Row 100
$ws.Range("H100").Value = 2929.9
$ws.Range("I100").Value = 955.5
$ws.Range("J100").Value = 5891.5
$ws.Range("K100").Value = 955.5
$ws.Range("L100").Value = 5891.5
$ws.Range("M100").Value = -414.5
$ws.Range("N100").Value = -6973.5

# Row 103
$ws.Range("H103").Value = 3796.8333
$ws.Range("I103").Value = 1947
$ws.Range("J103").Value = 4721.75
$ws.Range("K103").Value = 5841
$ws.Range("L103").Value = 14165.25
$ws.Range("M103").Value = -5255
$ws.Range("N103").Value = -15337.25

# Row 135
$ws.Range("H135").Value = 870
$ws.Range("I135").Value = 874.44446
$ws.Range("K135").Value = 7870.00014
$ws.Range("M135").Value = -5335.00014

# Row 136
$ws.Range("H136").Value = 30000
$ws.Range("I136").Value = 30000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 30000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -24900
$ws.Range("N136").ClearContents()

# Row 138
$ws.Range("H138").Value = 2606.4614
$ws.Range("J138").Value = 3717.348
$ws.Range("L138").Value = 11152.044
$ws.Range("N138").Value = -21432.044

# Row 139
$ws.Range("H139").Value = 69997.5
$ws.Range("J139").Value = 69997.5
$ws.Range("L139").Value = 69997.5
$ws.Range("N139").Value = -80277.5

# Row 141
$ws.Range("H141").Value = 6761.077
$ws.Range("I141").Value = 4649
$ws.Range("J141").Value = 8571.429
$ws.Range("K141").Value = 13947
$ws.Range("L141").Value = 25714.287
$ws.Range("M141").Value = -8767
$ws.Range("N141").Value = -36074.287

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 50002620
$ws.Range("I45").Value = 76924460
$ws.Range("J45").Value = 4915.143
$ws.Range("K45").Value = 76924460
$ws.Range("L45").Value = 4915.143
$ws.Range("M45").Value = -76924083
$ws.Range("N45").Value = -5669.143

# Row 61
$ws.Range("H61").Value = 8476.723
$ws.Range("I61").Value = 6324
$ws.Range("J61").Value = 12782.167
$ws.Range("K61").Value = 6324
$ws.Range("L61").Value = 12782.167
$ws.Range("M61").Value = -6112
$ws.Range("N61").Value = -13206.167

# Row 74
$ws.Range("H74").Value = 37043348
$ws.Range("I74").Value = 47625148
$ws.Range("J74").Value = 7050
$ws.Range("K74").Value = 47625148
$ws.Range("L74").Value = 7050
$ws.Range("M74").Value = -47624274
$ws.Range("N74").Value = -8798

# Row 77
$ws.Range("H77").Value = 37043348
$ws.Range("I77").Value = 47625148
$ws.Range("J77").Value = 7050
$ws.Range("K77").Value = 238125740
$ws.Range("L77").Value = 35250
$ws.Range("M77").Value = -238121372
$ws.Range("N77").Value = -43986

# Row 88
$ws.Range("H88").Value = 2332.8333
$ws.Range("J88").Value = 2211.25
$ws.Range("L88").Value = 2211.25
$ws.Range("N88").Value = -3023.25

# Row 91
$ws.Range("H91").Value = 2332.8333
$ws.Range("J91").Value = 2211.25
$ws.Range("L91").Value = 2211.25
$ws.Range("N91").Value = -5019.25

# Row 109
$ws.Range("H109").Value = 85964.336
$ws.Range("J109").Value = 85964.336
$ws.Range("L109").Value = 85964.336
$ws.Range("N109").Value = -88738.336

# Row 110
$ws.Range("H110").Value = 8270
$ws.Range("I110").Value = 6927.125
$ws.Range("K110").Value = 6927.125
$ws.Range("M110").Value = -4882.125

# Row 132
$ws.Range("H132").Value = 5539.294
$ws.Range("I132").Value = 2625.5386
$ws.Range("K132").Value = 7876.6158
$ws.Range("M132").Value = -5346.6158

# Row 136
$ws.Range("H136").Value = 8476.723
$ws.Range("I136").Value = 6324
$ws.Range("J136").Value = 12782.167
$ws.Range("K136").Value = 18972
$ws.Range("L136").Value = 38346.501
$ws.Range("M136").Value = -16422
$ws.Range("N136").Value = -43446.501

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2505
$ws.Range("I99").Value = 3129.75
$ws.Range("J99").Value = 1255.5
$ws.Range("K99").Value = 3129.75
$ws.Range("L99").Value = 1255.5
$ws.Range("M99").Value = -1631.75
$ws.Range("N99").Value = -4251.5

# Row 107
$ws.Range("H107").Value = 4814.2
$ws.Range("I107").Value = 4680.5
$ws.Range("J107").Value = 4903.3335
$ws.Range("K107").Value = 4680.5
$ws.Range("L107").Value = 4903.3335
$ws.Range("M107").Value = -2760.5
$ws.Range("N107").Value = -8743.333500000001

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 1649.3334
$ws.Range("I6").Value = 1099.2
$ws.Range("K6").Value = 1099.2
$ws.Range("M6").Value = -986.2

# Row 16
$ws.Range("H16").Value = 607.4167
$ws.Range("I16").Value = 661.125
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 661.125
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -374.125
$ws.Range("N16").Value = -1074

# Row 22
$ws.Range("H22").Value = 979.83636
$ws.Range("I22").Value = 741.67566
$ws.Range("K22").Value = 741.67566
$ws.Range("M22").Value = -391.67566

# Row 113
$ws.Range("H113").Value = 607.4167
$ws.Range("I113").Value = 661.125
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 661.125
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = 1508.875
$ws.Range("N113").Value = -4840

# Row 122
$ws.Range("H122").Value = 7470.2144
$ws.Range("I122").Value = 3846.5
$ws.Range("K122").Value = 11539.5
$ws.Range("M122").Value = -9089.5

# Row 132
$ws.Range("H132").Value = 5640.8096
$ws.Range("I132").Value = 4745.9287
$ws.Range("J132").Value = 7430.5713
$ws.Range("K132").Value = 14237.7861
$ws.Range("L132").Value = 22291.7139
$ws.Range("M132").Value = -11707.7861
$ws.Range("N132").Value = -27351.7139

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 12692.308
$ws.Range("J5").Value = 17446.555
$ws.Range("L5").Value = 52339.665
$ws.Range("N5").Value = -52563.665

# Row 32
$ws.Range("H32").Value = 18575
$ws.Range("J32").Value = 18575
$ws.Range("L32").Value = 55725
$ws.Range("N32").Value = -56291

# Row 135
$ws.Range("H135").Value = 12692.308
$ws.Range("J135").Value = 17446.555
$ws.Range("L135").Value = 157018.995
$ws.Range("N135").Value = -162088.995

# Row 136
$ws.Range("H136").Value = 1095.4445
$ws.Range("I136").Value = 1095.4445
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3286.3335
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 1813.6665
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5790.6
$ws.Range("I80").Value = 4322.222
$ws.Range("J80").Value = 19006
$ws.Range("K80").Value = 4322.222
$ws.Range("L80").Value = 19006
$ws.Range("M80").Value = -3324.222
$ws.Range("N80").Value = -21002

# Row 83
$ws.Range("H83").Value = 5790.6
$ws.Range("I83").Value = 4322.222
$ws.Range("J83").Value = 19006
$ws.Range("K83").Value = 21611.11
$ws.Range("L83").Value = 95030
$ws.Range("M83").Value = -16619.11
$ws.Range("N83").Value = -105014

# Row 102
$ws.Range("H102").Value = 5831.5
$ws.Range("I102").Value = 4910.375
$ws.Range("K102").Value = 4910.375
$ws.Range("M102").Value = -3288.375

# Row 122
$ws.Range("H122").Value = 5338.1333
$ws.Range("I122").Value = 4005.3333
$ws.Range("J122").Value = 10669.333
$ws.Range("K122").Value = 12015.9999
$ws.Range("L122").Value = 32007.999
$ws.Range("M122").Value = -9565.999899999999
$ws.Range("N122").Value = -36907.999

# Row 132
$ws.Range("H132").Value = 108988
$ws.Range("I132").Value = 147549
$ws.Range("K132").Value = 442647
$ws.Range("M132").Value = -440117

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 10299.667
$ws.Range("I40").Value = 7098.7144
$ws.Range("J40").Value = 14781
$ws.Range("K40").Value = 7098.7144
$ws.Range("L40").Value = 14781
$ws.Range("M40").Value = -6962.7144
$ws.Range("N40").Value = -15053

# Row 61
$ws.Range("H61").Value = 7862
$ws.Range("I61").Value = 2149.3333
$ws.Range("K61").Value = 2149.3333
$ws.Range("M61").Value = -1947.3333

# Row 100
$ws.Range("H100").Value = 3568.1724
$ws.Range("I100").Value = 2520.8262
$ws.Range("J100").Value = 7583
$ws.Range("K100").Value = 2520.8262
$ws.Range("L100").Value = 7583
$ws.Range("M100").Value = -1979.8262
$ws.Range("N100").Value = -8665

# Row 113
$ws.Range("H113").Value = 7862
$ws.Range("I113").Value = 2149.3333
$ws.Range("K113").Value = 2149.3333
$ws.Range("M113").Value = 20.66670000000022

# Row 122
$ws.Range("H122").Value = 7635
$ws.Range("I122").Value = 1950
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400

# Row 132
$ws.Range("H132").Value = 6284.3335
$ws.Range("I132").Value = 2533
$ws.Range("K132").Value = 7599
$ws.Range("M132").Value = -5069

# Row 136
$ws.Range("H136").Value = 5914.75
$ws.Range("I136").Value = 2911.4546
$ws.Range("J136").Value = 12522
$ws.Range("K136").Value = 8734.363799999999
$ws.Range("L136").Value = 37566
$ws.Range("M136").Value = -6184.363799999999
$ws.Range("N136").Value = -42666

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5877.8887
$ws.Range("I122").Value = 4237
$ws.Range("K122").Value = 12711
$ws.Range("M122").Value = -10261

# Row 136
$ws.Range("H136").Value = 4003.6667
$ws.Range("I136").Value = 2619.0667
$ws.Range("J136").Value = 8619
$ws.Range("K136").Value = 7857.2001
$ws.Range("L136").Value = 25857
$ws.Range("M136").Value = -5307.2001
$ws.Range("N136").Value = -30957

